$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (Price / Volume) to Text format so that
# numeric-looking strings (e.g. "29.324.50", "4.980") are not
# auto-converted to numbers by Excel's smart entry parsing.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.324.50"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.842.36"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "238.77"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").Value = "0.6305"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "0.07525"
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.2943"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "24.45"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.07694"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.826.39"
$ws.Range("E12").Value = "  -7.99%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.980"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.6784"
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.00001043"
$ws.Range("E15").Value = "  +5.20%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "82.98"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "2.089.28"
$ws.Range("E17").Value = "  -7.73%  "
$ws.Range("D18").Value = "6.142"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "29.367.67"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "228.30"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").Value = "0.9997"
$ws.Range("D23").Value = "7.426"
$ws.Range("E23").Value = "  -2.42%  "
$ws.Range("D25").Value = "156.51"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").Value = "0.1391"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "8.341"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("D29").Value = "1.455"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").Value = "1.270"
$ws.Range("E30").Value = "  +1.05%  "
$ws.Range("D31").Value = "0.05631"
$ws.Range("E31").Value = "  -3.23%  "
$ws.Range("D32").Value = "4.103"
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("D34").Value = "1.829"
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("D36").Value = "0.7087"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "1.245.34"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").Value = "6.239"
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("D42").Value = "0.9025"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").Value = "0.9993"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "101.82"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").Value = "65.65"
$ws.Range("E45").Value = "  -2.40%  "
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("D47").Value = "7.109"
$ws.Range("E47").Value = "  -2.74%  "
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("D50").Value = "8.898"
$ws.Range("E50").Value = "  -3.08%  "
$ws.Range("E51").Value = "  -0.30%  "
